$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header labels: "Item ID" -> "ItemID", "Item Name" -> "ItemNameE"
$ws.Range("A1").Value = "ItemID"
$ws.Range("B1").Value = "ItemNameE"

# Remove the border formatting from the header cells (they now use the
# default style instead of the bordered one)
$ws.Range("A1:B1").Borders.LineStyle = -4142

# Move the active selection to G6
$ws.Range("G6").Select()
